$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Report Generated On" timestamp
$ws.Range("D5").Value = "Report Generated On: 08/26/2025 10:01 AM"

# Update Total Billed Amount
$ws.Range("C8").Value = 34.4

# Clear the Scope ID # value (was "#NO MATCH", now blank)
$ws.Range("G10").Value = ""

# Update per-line-item pricing and TOTAL pricing
$ws.Range("H16").Value = 34.4
$ws.Range("H17").Value = 34.4
